# "Final touches for C23A"
# The workbook is a crossword/"Clue" grid. Each clue-start cell is rendered
# with a highlighted fill (style index 7 / 8) and a shared string that has
# a trailing "#" marker (e.g. "B#" instead of "B"). This edit moves the
# clue-start marker for six clues to the correct cell within their row,
# leaving the underlying letter unchanged.
#
# For each move below we:
#   1. Copy the *formats* (fill + underlying style) from the cell that
#      currently carries the marker onto the cell that should receive it,
#      then set that destination cell's value to the "#" variant string.
#   2. Restore the original (no-marker) cell back to the plain style by
#      copying formats from a neighboring plain cell in the same row, then
#      set its value back to the plain (no "#") string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Move-ClueMarker {
    param(
        [string]$From,
        [string]$To,
        [string]$PlainDonor,
        [string]$MarkedText,
        [string]$PlainText
    )

    # Give the destination cell the highlighted "marker" formatting.
    $ws.Range($From).Copy() | Out-Null
    $ws.Range($To).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($To).Value = $MarkedText

    # Restore the source cell to plain formatting.
    $ws.Range($PlainDonor).Copy() | Out-Null
    $ws.Range($From).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws.Range($From).Value = $PlainText

    $excel.CutCopyMode = 0
}

# Clue "B": marker moves from D3 to E4
Move-ClueMarker "D3" "E4" "C3" "B#" "B"

# Clue "S": marker moves from M4 to L8
Move-ClueMarker "M4" "L8" "L4" "S#" "S"

# Clue "G": marker moves from D15 to J17
Move-ClueMarker "D15" "J17" "C15" "G#" "G"

# Clue "L": marker moves from AC21 to AC20
Move-ClueMarker "AC21" "AC20" "AB21" "L#" "L"

# Clue "H": marker moves from D22 to D24
Move-ClueMarker "D22" "D24" "C22" "H#" "H"

# Clue "D": marker moves from O22 to O24
Move-ClueMarker "O22" "O24" "N22" "D#" "D"

# Reflect the cursor/selection position that was active when the file was
# last saved.
$ws.Range("L10").Select() | Out-Null
